# Updates cryptos list values (price/volume) scraped for this commit.
# Cells are written as Text (matching the source inlineStr cells); numeric-looking
# strings are prefixed with a leading apostrophe so Excel keeps them as text instead
# of auto-converting to a Number (which would also drop formatting, e.g. "38.90" -> 38.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.569.63"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "3.187.19"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'602.51"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").Value = "'155.63"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.185.74"
$ws.Range("E8").Value = "  -1.34%  "

$ws.Range("E9").Value = "  +2.35%  "

$ws.Range("E10").Value = "  -2.33%  "

$ws.Range("D11").Value = "'5.67"
$ws.Range("E11").Value = "  -8.67%  "

$ws.Range("D12").Value = "'0.515"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  -1.68%  "

$ws.Range("D14").Value = "'38.90"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "3.711.68"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").Value = "66.573.65"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").Value = "'7.46"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "3.182.80"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").Value = "'514.56"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").Value = "'15.49"
$ws.Range("E21").Value = "  -2.18%  "

$ws.Range("D22").Value = "'0.736"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").Value = "'8.16"
$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("D24").Value = "'14.96"
$ws.Range("E24").Value = "  -3.00%  "

$ws.Range("D25").Value = "'84.87"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27/28: RenderToken and PancakeSwap swapped rank order
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'9.30"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'3.01"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("E29").Value = "  +6.43%  "

$ws.Range("D30").Value = "'3.12"
$ws.Range("E30").Value = "  +7.09%  "

$ws.Range("D31").Value = "'7.00"
$ws.Range("E31").Value = "  +1.31%  "

$ws.Range("D32").Value = "'28.19"
$ws.Range("E32").Value = "  -0.86%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("E34").Value = "  -1.40%  "

$ws.Range("E35").Value = "  -2.13%  "

$ws.Range("D36").Value = "'516.30"
$ws.Range("E36").Value = "  +4.70%  "

$ws.Range("D37").Value = "'54.84"
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("E38").Value = "  -3.47%  "

$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "'0.128"
$ws.Range("E40").Value = "  +6.45%  "

$ws.Range("D41").Value = "'8.89"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").Value = "0.0₃0687"
$ws.Range("E42").Value = "  +4.93%  "

$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("E44").Value = "  -8.17%  "

$ws.Range("D45").Value = "'2.47"
$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("D46").Value = "2.855.96"
$ws.Range("E46").Value = "  -6.28%  "

$ws.Range("D47").Value = "'28.31"
$ws.Range("E47").Value = "  -3.60%  "

$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +2.27%  "

$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").Value = "'2.59"
$ws.Range("E51").Value = "  +5.95%  "
